$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Split "BaseX aN514_Tomme_mapper.xq ..." run into several runs
#    separated by <w:proofErr> spell/grammar-check markers, and
#    re-split the filename text between "aN514_Tomme_" and
#    "mapper.xq".
# -----------------------------------------------------------------
$findRng = $d.Content
$found = $findRng.Find.Execute("BaseX", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $findRng.Expand(4)  # wdParagraph - grab the whole paragraph incl. its mark

    $xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>BaseX</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:t>aN514_Tomme_</w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:t>mapper.xq</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:t xml:space="preserve"> for &#229; sjekke antall tomme mapper som ikke utg&#229;r.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
    $findRng.InsertXML($xml)
}

# -----------------------------------------------------------------
# 2. Remove the whole "distribution" table and leave a single blank
#    paragraph in its place (the paragraph that already followed
#    the table is left untouched).
# -----------------------------------------------------------------
if ($d.Tables.Count -ge 1) {
    $t = $d.Tables.Item(1)
    $tStart = $t.Range.Start
    $t.Delete()

    $gap = $d.Range($tStart - 1, $tStart - 1)
    $gap.Text = [string][char]13
}
